# Update the public EPEX Spot prices workbook:
# 1. "Prix Spot" sheet: insert a new date column (01-nov) before the
#    01-oct. column, shifting the October/beyond columns one position to
#    the right, and fill the new column with "-" placeholders (no data
#    yet for 01-nov).
# 2. "Gaz" sheet: append the new day's row (2025-10-30).
# 3. "CO2" sheet: append the new day's row (2025-10-30).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "Prix Spot": insert column before DF (i.e. before 01-oct.)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Insert a new column at DF (column 110); this shifts old DF:EJ -> DG:EK
# automatically, carrying over values, types and styles.
$ws1.Range("DF1").EntireColumn.Insert()

# New header cell for the inserted column.
$ws1.Range("DF1").Value = "01-nov"

# There is no data yet for 01-nov, so every data row gets a "-" placeholder,
# matching the convention used for other future/missing days.
for ($row = 2; $row -le 25; $row++) {
    $ws1.Cells.Item($row, 110).Value = "-"
}

# ---------------------------------------------------------------------
# 2) Sheet "Gaz": append row 138 (2025-10-30)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Gaz")

# Use a leading apostrophe so the ISO-formatted date string is kept as
# text (matching every other Date cell in this column) instead of being
# auto-converted to a date serial number, then drop the resulting
# "quote prefix" formatting so the cell stays visually identical to its
# neighbours.
$ws2.Range("A138").Value = "'2025-10-30"
$ws2.Range("A138").ClearFormats()
$ws2.Range("B138").Value = 29.8

# ---------------------------------------------------------------------
# 3) Sheet "CO2": append row 138 (2025-10-30)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Range("A138").Value = "'2025-10-30"
$ws3.Range("A138").ClearFormats()
$ws3.Range("B138").Value = 78.36
